# Insert a new daily price record as row 499 ("Fruta / hortaliza, semanal"),
# pushing the existing rows 499-560 down to 500-561.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 499; Excel shifts row 499..560 down to 500..561
# and the new row 499 inherits formatting (incl. date style) from the row above.
$ws.Rows.Item(499).EntireRow.Insert()

# Populate the new row 499 with the new record's values.
$ws.Range("A499").Value = 8
$ws.Range("B499").Value = "Terminal La Palmera de La Serena"
$ws.Range("C499").Value = "Coquimbo"
$ws.Range("D499").Value = 45142
$ws.Range("E499").Value = 4
$ws.Range("F499").Value = 100112032
$ws.Range("G499").Value = "Zapallo italiano"
$ws.Range("H499").Value = "Sin especificar"
$ws.Range("I499").Value = "Primera"
$ws.Range("J499").Value = 400
$ws.Range("K499").Value = 12000
$ws.Range("L499").Value = 13000
$ws.Range("M499").Value = 12500
$ws.Range("N499").Value = "`$/caja 50 unidades"
$ws.Range("O499").Value = "Región de Arica y Parinacota"
$ws.Range("P499").Value = 250
$ws.Range("Q499").Value = 50
$ws.Range("R499").Value = "Hortaliza"
